# CKU overall advancement. Some minor code changes
#
# The "empires" sheet already has a "government" header in column E (E1).
# This adds the corresponding value for the "iberia" row (row 2):
# "feudal_government".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "feudal_government"

# Excel auto-sizes a freshly populated column to fit its new content.
$ws.Columns.Item(5).AutoFit() | Out-Null

# Leave the cursor on the cell that was just edited, matching the
# author's final selection when the file was saved.
$ws.Range("E2").Select() | Out-Null
